# Update the MauiTemplates readme.docx footer copyright year range and
# turn on "Different odd & even" headers/footers (the document already has
# "Different first page" enabled via titlePg). Enabling odd/even causes
# Word to mint separate primary/even/first header & footer parts.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- 1. Update the copyright year text in the (primary/default) footer ---
$footerDefault = $sec.Footers(1)
$footerDefault.Range.Find.Execute(
    " 2022 Vijay Anand E G", $true, $false, $false, $false, $false,
    $true, 1, $false, " 2021-2024 Vijay Anand E G", 2)

# --- 2. Turn on distinct odd/even page headers & footers ---
# This, combined with the already-enabled "different first page" setting,
# makes Word create 3 header parts (even / default / first) and 3 footer
# parts (even / default / first) and wire up the section's references.
$sec.PageSetup.OddAndEvenPagesHeaderFooter = -1

# Touch the even & first-page headers/footers so Word actually mints the
# (empty) parts for them rather than leaving the section pointing only at
# the primary header/footer.
$sec.Headers(2).Range.Text = $sec.Headers(2).Range.Text
$sec.Headers(3).Range.Text = $sec.Headers(3).Range.Text
$sec.Footers(2).Range.Text = $sec.Footers(2).Range.Text
$sec.Footers(3).Range.Text = $sec.Footers(3).Range.Text

Write-Output ("Default footer: " + $footerDefault.Range.Text)
